$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the last data row's formatting, then insert it before row 29 (the total row),
# which shifts the total row (and everything below) down and keeps formatting intact.
$ws.Rows("28:28").Copy()
$ws.Rows("29:29").Insert()

# Populate the new row 29 with the next schedule entry
$ws.Range("B29").Value = 23
$ws.Range("C29").Value = "Revision-11_File duplicate checks updated"
$ws.Range("D29").Value = "16 - 01 - 2020"
$ws.Range("E29").Value = "Done"
$ws.Range("F29").Value = 1
